$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.926.39"
$ws.Range("E2").Value = "  -0.69%  "

# Row 3
$ws.Range("D3").Value = "2.578.69"
$ws.Range("E3").Value = "  +0.47%  "

# Row 4
$ws.Range("D4").Value = "'1.00"

# Row 5
$ws.Range("D5").Value = "'582.08"
$ws.Range("E5").Value = "  -0.46%  "

# Row 6
$ws.Range("D6").Value = "'144.38"
$ws.Range("E6").Value = "  -2.60%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  -2.08%  "

# Row 9
$ws.Range("D9").Value = "'0.107"
$ws.Range("E9").Value = "  -2.17%  "

# Row 10
$ws.Range("E10").Value = "  -1.60%  "

# Row 11
$ws.Range("E11").Value = "  -0.75%  "

# Row 12
$ws.Range("E12").Value = "  -2.57%  "

# Row 13
$ws.Range("D13").Value = "'27.02"
$ws.Range("E13").Value = "  -1.98%  "

# Row 14
$ws.Range("D14").Value = "3.041.67"
$ws.Range("E14").Value = "  +0.47%  "

# Row 15
$ws.Range("D15").Value = "62.814.93"
$ws.Range("E15").Value = "  -0.71%  "

# Row 16
$ws.Range("D16").Value = "'0.0000145"
$ws.Range("E16").Value = "  -2.71%  "

# Row 17
$ws.Range("D17").Value = "2.578.78"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("D18").Value = "'11.09"
$ws.Range("E18").Value = "  -2.61%  "

# Row 19
$ws.Range("D19").Value = "'340.53"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20
$ws.Range("D20").Value = "'4.33"
$ws.Range("E20").Value = "  -2.18%  "

# Row 21
$ws.Range("D21").Value = "'6.63"
$ws.Range("E21").Value = "  -3.49%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").Value = "'67.36"
$ws.Range("E23").Value = "  +0.75%  "

# Row 24
$ws.Range("D24").Value = "'1.59"
$ws.Range("E24").Value = "  +7.16%  "

# Row 25
$ws.Range("E25").Value = "  -2.39%  "

# Row 26
$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "  -3.46%  "

# Row 27
$ws.Range("D27").Value = "'8.00"
$ws.Range("E27").Value = "  -2.52%  "

# Row 28
$ws.Range("E28").Value = "  +0.04%  "

# Row 29
$ws.Range("D29").Value = "'8.25"
$ws.Range("E29").Value = "  -3.71%  "

# Row 30
$ws.Range("E30").Value = "  -3.57%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'460.47"
$ws.Range("E31").Value = "  -1.58%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0798"
$ws.Range("E32").Value = "  -3.76%  "

# Row 33
$ws.Range("D33").Value = "'1.65"
$ws.Range("E33").Value = "  +1.32%  "

# Row 34
$ws.Range("D34").Value = "'176.10"
$ws.Range("E34").Value = "  -0.39%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").Value = "'0.396"
$ws.Range("E36").Value = "  -2.62%  "

# Row 37
$ws.Range("D37").Value = "'18.86"
$ws.Range("E37").Value = "  -2.32%  "

# Row 38
$ws.Range("D38").Value = "'4.47"
$ws.Range("E38").Value = "  -1.40%  "

# Row 39
$ws.Range("E39").Value = "  -0.02%  "

# Row 40
$ws.Range("E40").Value = "  -3.69%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'158.20"
$ws.Range("E41").Value = "  +4.20%  "

# Row 42
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'39.98"
$ws.Range("E42").Value = "  +0.68%  "

# Row 43
$ws.Range("E43").Value = "  -3.83%  "

# Row 44
$ws.Range("D44").Value = "'0.637"
$ws.Range("E44").Value = "  +3.29%  "

# Row 45
$ws.Range("D45").Value = "'21.14"
$ws.Range("E45").Value = "  -0.36%  "

# Row 46
$ws.Range("D46").Value = "'0.0537"
$ws.Range("E46").Value = "  -3.08%  "

# Row 47
$ws.Range("E47").Value = "  -2.44%  "

# Row 48
$ws.Range("E48").Value = "  -2.48%  "

# Row 49
$ws.Range("D49").Value = "'17.96"
$ws.Range("E49").Value = "  -2.99%  "

# Row 50
$ws.Range("D50").Value = "'11.41"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("E51").Value = "  -4.67%  "
